$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 155.2
$ws.Range("I55").Value = 156
$ws.Range("J55").Value = 153.33333
$ws.Range("K55").Value = 156
$ws.Range("L55").Value = 153.33333
$ws.Range("M55").Value = 58
$ws.Range("N55").Value = -581.3333299999999
$ws.Range("H100").Value = 3491.5
$ws.Range("I100").Value = 3364.375
$ws.Range("K100").Value = 3364.375
$ws.Range("M100").Value = -2823.375
$ws.Range("H116").Value = 3399
$ws.Range("I116").Value = 3498.75
$ws.Range("K116").Value = 3498.75
$ws.Range("M116").Value = -56.75
$ws.Range("H131").Value = 71429290
$ws.Range("I131").Value = 71429290
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 214287870
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -214282830
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 3000
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H10").Value = 24595.455
$ws.Range("J10").Value = 24595.455
$ws.Range("L10").Value = 24595.455
$ws.Range("N10").Value = -24935.455
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H21").Value = 12666.667
$ws.Range("I21").Value = 1500
$ws.Range("J21").Value = 14900
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 14900
$ws.Range("M21").Value = -1126
$ws.Range("N21").Value = -15648
$ws.Range("H22").Value = 18000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H26").Value = 8398
$ws.Range("I26").Value = 1007
$ws.Range("J26").Value = 15789
$ws.Range("K26").Value = 1007
$ws.Range("L26").Value = 15789
$ws.Range("M26").Value = -677
$ws.Range("N26").Value = -16449
$ws.Range("H27").Value = 11000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 11000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 11000
$ws.Range("N27").Value = -11368
$ws.Range("M27").ClearContents()
$ws.Range("H33").Value = 9788
$ws.Range("I33").Value = 9676
$ws.Range("J33").Value = 9900
$ws.Range("K33").Value = 9676
$ws.Range("L33").Value = 9900
$ws.Range("M33").Value = -9347
$ws.Range("N33").Value = -10558
$ws.Range("H34").Value = 9620
$ws.Range("I34").Value = 10000
$ws.Range("J34").Value = 9525
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 9525
$ws.Range("M34").Value = -9729
$ws.Range("N34").Value = -10067
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H39").Value = 18325
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 18325
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 18325
$ws.Range("N39").Value = -19365
$ws.Range("M39").ClearContents()
$ws.Range("H40").Value = 13583.333
$ws.Range("J40").Value = 13583.333
$ws.Range("L40").Value = 13583.333
$ws.Range("N40").Value = -13935.333
$ws.Range("H51").Value = 10350
$ws.Range("J51").Value = 10350
$ws.Range("L51").Value = 10350
$ws.Range("N51").Value = -11862
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H102").Value = 2406.111
$ws.Range("I102").Value = 1967.5
$ws.Range("J102").Value = 3283.3333
$ws.Range("K102").Value = 1967.5
$ws.Range("L102").Value = 3283.3333
$ws.Range("M102").Value = -345.5
$ws.Range("N102").Value = -6527.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 990
$ws.Range("I8").Value = 990
$ws.Range("K8").Value = 990
$ws.Range("M8").Value = -850
$ws.Range("H10").Value = 2000.6666
$ws.Range("I10").Value = 2501.3333
$ws.Range("J10").Value = 1500
$ws.Range("K10").Value = 2501.3333
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = -2361.3333
$ws.Range("N10").Value = -1780
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20454
$ws.Range("H16").Value = 40000
$ws.Range("J16").Value = 40000
$ws.Range("L16").Value = 40000
$ws.Range("N16").Value = -40340
$ws.Range("H22").Value = 306
$ws.Range("I22").Value = 360.83334
$ws.Range("J22").Value = 240.2
$ws.Range("K22").Value = 360.83334
$ws.Range("L22").Value = 240.2
$ws.Range("M22").Value = -187.83334
$ws.Range("N22").Value = -586.2
$ws.Range("H32").Value = 26000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 26000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 26000
$ws.Range("N32").Value = -26768
$ws.Range("M32").ClearContents()
$ws.Range("H33").Value = 23142.834
$ws.Range("I33").Value = 4180.6665
$ws.Range("J33").Value = 42105
$ws.Range("K33").Value = 4180.6665
$ws.Range("L33").Value = 42105
$ws.Range("M33").Value = -3844.6665
$ws.Range("N33").Value = -42777
$ws.Range("H82").Value = 13799.857
$ws.Range("I82").Value = 6649.75
$ws.Range("J82").Value = 23333.334
$ws.Range("K82").Value = 6649.75
$ws.Range("L82").Value = 23333.334
$ws.Range("M82").Value = -6266.75
$ws.Range("N82").Value = -24099.334
$ws.Range("H85").Value = 13799.857
$ws.Range("I85").Value = 6649.75
$ws.Range("J85").Value = 23333.334
$ws.Range("K85").Value = 6649.75
$ws.Range("L85").Value = 23333.334
$ws.Range("M85").Value = -5323.75
$ws.Range("N85").Value = -25985.334
$ws.Range("H103").Value = 35000.332
$ws.Range("J103").Value = 35000.332
$ws.Range("L103").Value = 35000.332
$ws.Range("N103").Value = -37344.332
$ws.Range("H105").Value = 2398.2292
$ws.Range("I105").Value = 2378.0264
$ws.Range("J105").Value = 2475
$ws.Range("K105").Value = 2378.0264
$ws.Range("L105").Value = 2475
$ws.Range("M105").Value = -631.0264000000002
$ws.Range("N105").Value = -5969

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 60000
$ws.Range("J63").Value = 60000
$ws.Range("L63").Value = 60000
$ws.Range("N63").Value = -61372
$ws.Range("H66").Value = 60000
$ws.Range("J66").Value = 60000
$ws.Range("L66").Value = 180000
$ws.Range("N66").Value = -186864
$ws.Range("H80").Value = 22333.334
$ws.Range("J80").Value = 22333.334
$ws.Range("L80").Value = 22333.334
$ws.Range("N80").Value = -24579.334
$ws.Range("H83").Value = 22333.334
$ws.Range("J83").Value = 22333.334
$ws.Range("L83").Value = 67000.00199999999
$ws.Range("N83").Value = -78232.00199999999
$ws.Range("H132").Value = 825107.5
$ws.Range("I132").Value = 2104.6775
$ws.Range("J132").Value = 2647470.8
$ws.Range("K132").Value = 6314.032499999999
$ws.Range("L132").Value = 7942412.399999999
$ws.Range("M132").Value = -3784.032499999999
$ws.Range("N132").Value = -7947472.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 311
$ws.Range("I40").Value = 59.5
$ws.Range("J40").Value = 646.3333
$ws.Range("K40").Value = 238
$ws.Range("L40").Value = 2585.3332
$ws.Range("M40").Value = -169
$ws.Range("N40").Value = -2723.3332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 11499.5
$ws.Range("J15").Value = 11499.5
$ws.Range("L15").Value = 11499.5
$ws.Range("N15").Value = -12075.5
$ws.Range("H81").Value = 11499.5
$ws.Range("J81").Value = 11499.5
$ws.Range("L81").Value = 11499.5
$ws.Range("N81").Value = -13495.5
$ws.Range("H84").Value = 11499.5
$ws.Range("J84").Value = 11499.5
$ws.Range("L84").Value = 34498.5
$ws.Range("N84").Value = -44482.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3649.818
$ws.Range("I7").Value = 2884.2856
$ws.Range("J7").Value = 4989.5
$ws.Range("K7").Value = 2884.2856
$ws.Range("L7").Value = 4989.5
$ws.Range("M7").Value = -2772.2856
$ws.Range("N7").Value = -5213.5
$ws.Range("H40").Value = 1521.56
$ws.Range("I40").Value = 1448.7646
$ws.Range("J40").Value = 1676.25
$ws.Range("K40").Value = 1448.7646
$ws.Range("L40").Value = 1676.25
$ws.Range("M40").Value = -1312.7646
$ws.Range("N40").Value = -1948.25
$ws.Range("H126").Value = 3649.818
$ws.Range("I126").Value = 2884.2856
$ws.Range("J126").Value = 4989.5
$ws.Range("K126").Value = 8652.856800000001
$ws.Range("L126").Value = 14968.5
$ws.Range("M126").Value = -6182.856800000001
$ws.Range("N126").Value = -19908.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1094.5
$ws.Range("I126").Value = 970.6667
$ws.Range("J126").Value = 1218.3334
$ws.Range("K126").Value = 2912.0001
$ws.Range("L126").Value = 3655.0002
$ws.Range("M126").Value = -442.0001000000002
$ws.Range("N126").Value = -8595.0002
